# Daily refresh of the crypto symbol list (Price / Volume(1h) / Hora columns,
# plus a few re-ranked Coin/Link pairs) as produced by the GitHub Actions job.
#
# Price, Volume(1h) and Hora are stored as literal TEXT in this sheet (e.g. a
# trailing-zero price like "27.40" or a literal "4.50%" string), not as actual
# numbers/percentages. Writing a bare numeric-looking string via COM would make
# Excel auto-convert the cell to a Number, silently dropping formatting such as
# trailing zeros or the '%'. Prefixing the literal with a leading apostrophe
# (the same trick used when typing into the Excel UI) forces it to stay text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '''257.69'
$ws.Cells.Item(2, 5).Value = '''4.50%'
$ws.Cells.Item(2, 7).Value = '''2'
# Row 3
$ws.Cells.Item(3, 4).Value = '''27.29'
$ws.Cells.Item(3, 5).Value = '''-4.98%'
$ws.Cells.Item(3, 7).Value = '''2'
# Row 4
$ws.Cells.Item(4, 4).Value = '''5.214'
$ws.Cells.Item(4, 5).Value = '''-1.23%'
$ws.Cells.Item(4, 7).Value = '''2'
# Row 5
$ws.Cells.Item(5, 4).Value = '''0.05939'
$ws.Cells.Item(5, 5).Value = '''3.48%'
$ws.Cells.Item(5, 7).Value = '''2'
# Row 6
$ws.Cells.Item(6, 4).Value = '''6.709'
$ws.Cells.Item(6, 5).Value = '''0.61%'
$ws.Cells.Item(6, 7).Value = '''2'
# Row 7
$ws.Cells.Item(7, 4).Value = '''0.8689'
$ws.Cells.Item(7, 5).Value = '''0.86%'
$ws.Cells.Item(7, 7).Value = '''2'
# Row 8
$ws.Cells.Item(8, 4).Value = '''1.001'
$ws.Cells.Item(8, 5).Value = '''8.60%'
$ws.Cells.Item(8, 7).Value = '''2'
# Row 9
$ws.Cells.Item(9, 2).Value = 'WazirX'
$ws.Cells.Item(9, 3).Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Cells.Item(9, 4).Value = '''0.1417'
$ws.Cells.Item(9, 5).Value = '''1.49%'
$ws.Cells.Item(9, 7).Value = '''2'
# Row 10
$ws.Cells.Item(10, 2).Value = 'MandalaExchangeToken'
$ws.Cells.Item(10, 3).Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Cells.Item(10, 4).Value = '''0.07194'
$ws.Cells.Item(10, 5).Value = '''0.60%'
$ws.Cells.Item(10, 7).Value = '''2'
# Row 11
$ws.Cells.Item(11, 2).Value = 'BitrueCoin'
$ws.Cells.Item(11, 3).Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Cells.Item(11, 4).Value = '''0.03149'
$ws.Cells.Item(11, 5).Value = '''0.66%'
$ws.Cells.Item(11, 7).Value = '''2'
# Row 12
$ws.Cells.Item(12, 2).Value = 'BitMartToken'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Cells.Item(12, 4).Value = '''0.09247'
$ws.Cells.Item(12, 5).Value = '''0.14%'
$ws.Cells.Item(12, 7).Value = '''2'
# Row 13
$ws.Cells.Item(13, 2).Value = 'BitForexToken'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Cells.Item(13, 4).Value = '''0.001547'
$ws.Cells.Item(13, 5).Value = '''1.15%'
$ws.Cells.Item(13, 7).Value = '''2'
# Row 14
$ws.Cells.Item(14, 2).Value = 'One'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Cells.Item(14, 4).Value = '''0.0006070'
$ws.Cells.Item(14, 5).Value = '''1.12%'
$ws.Cells.Item(14, 7).Value = '''2'
# Row 15
$ws.Cells.Item(15, 4).Value = '''0.006144'
$ws.Cells.Item(15, 5).Value = '''4.49%'
$ws.Cells.Item(15, 7).Value = '''2'
# Row 16
$ws.Cells.Item(16, 4).Value = '''3.504'
$ws.Cells.Item(16, 5).Value = '''-0.01%'
$ws.Cells.Item(16, 7).Value = '''2'
# Row 17
$ws.Cells.Item(17, 4).Value = '''3.269'
$ws.Cells.Item(17, 5).Value = '''1.08%'
$ws.Cells.Item(17, 7).Value = '''2'
# Row 18
$ws.Cells.Item(18, 4).Value = '''2.206'
$ws.Cells.Item(18, 5).Value = '''-1.38%'
$ws.Cells.Item(18, 7).Value = '''2'
# Row 19
$ws.Cells.Item(19, 4).Value = '''0.3146'
$ws.Cells.Item(19, 5).Value = '''0.56%'
$ws.Cells.Item(19, 7).Value = '''2'
# Row 20
$ws.Cells.Item(20, 4).Value = '''0.03543'
$ws.Cells.Item(20, 5).Value = '''5.35%'
$ws.Cells.Item(20, 7).Value = '''2'
# Row 21
$ws.Cells.Item(21, 4).Value = '''0.1307'
$ws.Cells.Item(21, 5).Value = '''-0.40%'
$ws.Cells.Item(21, 7).Value = '''2'
# Row 22
$ws.Cells.Item(22, 4).Value = '''3.584'
$ws.Cells.Item(22, 5).Value = '''1.34%'
$ws.Cells.Item(22, 7).Value = '''2'
# Row 23
$ws.Cells.Item(23, 4).Value = '''0.04262'
$ws.Cells.Item(23, 5).Value = '''2.42%'
$ws.Cells.Item(23, 7).Value = '''2'
# Row 24
$ws.Cells.Item(24, 4).Value = '''0.1349'
$ws.Cells.Item(24, 5).Value = '''-2.03%'
$ws.Cells.Item(24, 7).Value = '''2'
# Row 25
$ws.Cells.Item(25, 4).Value = '''0.001219'
$ws.Cells.Item(25, 5).Value = '''-0.29%'
$ws.Cells.Item(25, 7).Value = '''2'
# Row 26
$ws.Cells.Item(26, 4).Value = '''0.004518'
$ws.Cells.Item(26, 5).Value = '''-10.31%'
$ws.Cells.Item(26, 7).Value = '''2'
# Row 27
$ws.Cells.Item(27, 5).Value = '''0.06%'
$ws.Cells.Item(27, 7).Value = '''2'
# Row 28
$ws.Cells.Item(28, 5).Value = '''-22.97%'
$ws.Cells.Item(28, 7).Value = '''2'
# Row 29
$ws.Cells.Item(29, 7).Value = '''2'
# Row 30
$ws.Cells.Item(30, 7).Value = '''2'
# Row 31
$ws.Cells.Item(31, 7).Value = '''2'
# Row 32
$ws.Cells.Item(32, 7).Value = '''2'
# Row 33
$ws.Cells.Item(33, 7).Value = '''2'
# Row 34
$ws.Cells.Item(34, 7).Value = '''2'
# Row 35
$ws.Cells.Item(35, 7).Value = '''2'
# Row 36
$ws.Cells.Item(36, 7).Value = '''2'
# Row 37
$ws.Cells.Item(37, 7).Value = '''2'
# Row 38
$ws.Cells.Item(38, 7).Value = '''2'
# Row 39
$ws.Cells.Item(39, 7).Value = '''2'
# Row 40
$ws.Cells.Item(40, 5).Value = '''-0.24%'
$ws.Cells.Item(40, 7).Value = '''2'
# Row 41
$ws.Cells.Item(41, 2).Value = 'BKEXToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Cells.Item(41, 4).Value = '''0.1105'
$ws.Cells.Item(41, 5).Value = '''2.32%'
$ws.Cells.Item(41, 7).Value = '''2'
# Row 42
$ws.Cells.Item(42, 2).Value = 'KickToken'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Cells.Item(42, 4).Value = '''0.003989'
$ws.Cells.Item(42, 5).Value = '''-29.68%'
$ws.Cells.Item(42, 7).Value = '''2'
# Row 43
$ws.Cells.Item(43, 4).Value = '''0.002341'
$ws.Cells.Item(43, 5).Value = '''-3.25%'
$ws.Cells.Item(43, 7).Value = '''2'
# Row 44
$ws.Cells.Item(44, 4).Value = '''0.01049'
$ws.Cells.Item(44, 5).Value = '''9.84%'
$ws.Cells.Item(44, 7).Value = '''2'
# Row 45
$ws.Cells.Item(45, 4).Value = '''0.00005488'
$ws.Cells.Item(45, 5).Value = '''3.99%'
$ws.Cells.Item(45, 7).Value = '''2'
# Row 46
$ws.Cells.Item(46, 5).Value = '''0.12%'
$ws.Cells.Item(46, 7).Value = '''2'
# Row 47
$ws.Cells.Item(47, 5).Value = '''28.50%'
$ws.Cells.Item(47, 7).Value = '''2'
# Row 48
$ws.Cells.Item(48, 4).Value = '''0.002225'
$ws.Cells.Item(48, 5).Value = '''2.26%'
$ws.Cells.Item(48, 7).Value = '''2'
# Row 49
$ws.Cells.Item(49, 5).Value = '''0.12%'
$ws.Cells.Item(49, 7).Value = '''2'
# Row 50
$ws.Cells.Item(50, 5).Value = '''0.12%'
$ws.Cells.Item(50, 7).Value = '''2'
# Row 51
$ws.Cells.Item(51, 7).Value = '''2'
